# comento descarga de archivos en archivosdiarios.php
#
# The underlying report (reporte_transacciones) was regenerated: the
# "Info adicional comercio" column now carries a customer code, and the
# "Fecha trx" timestamps for the most recent batch of transactions were
# normalized to 2025-10-15 (keeping the original time-of-day/offset).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Info adicional comercio" (column E) -> customer code, for every data row.
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 5).Value = "C11672"
}

# "Fecha trx" (column K) -> re-dated to 2025-10-15, same time-of-day.
$ws.Cells.Item(2, 11).Value  = "2025-10-15T15:12:07-03:00"
$ws.Cells.Item(3, 11).Value  = "2025-10-15T15:13:07-03:00"
$ws.Cells.Item(4, 11).Value  = "2025-10-15T15:13:51-03:00"
$ws.Cells.Item(5, 11).Value  = "2025-10-15T15:15:17-03:00"
$ws.Cells.Item(6, 11).Value  = "2025-10-15T15:15:58-03:00"
$ws.Cells.Item(7, 11).Value  = "2025-10-15T17:11:31-03:00"
$ws.Cells.Item(8, 11).Value  = "2025-10-15T17:12:23-03:00"
$ws.Cells.Item(9, 11).Value  = "2025-10-15T17:14:20-03:00"
$ws.Cells.Item(10, 11).Value = "2025-10-15T17:32:08-03:00"
$ws.Cells.Item(11, 11).Value = "2025-10-15T17:32:45-03:00"
